# Nesto_TestCases.xlsx - "eighteenth commit with customer side test completed"
# Adds two new customer test cases (TC_CUST_05 "Edit Customer Details" and
# TC_CUST_06 "Delete Customer") to the Customer_Tests sheet, then leaves that
# sheet as the active/selected tab (mirroring the author's last editing
# position in Excel before save).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer_Tests")

# TC_CUST_05 - Edit Customer Details
$ws.Range("A13").Value = "TC_CUST_05"
$ws.Range("B13").Value = "Edit Customer Details"
$ws.Range("C13").Value = "1.Click ""Edit Button"" at ""//tr[td[text()='auto@nesto.com']]//a[contains(@class, 'btn-outline-primary')]"""
$ws.Range("C14").Value = "2.Type ""Automation User Updated"" at ""//input[@name='name']"""
$ws.Range("C15").Value = "3.Click ""Save Customer"" at ""//button[contains(@class, 'btn-save')]"""
$ws.Range("C16").Value = "4.Verify text ""{DB_QUERY}SELECT name FROM customers WHERE email='auto@nesto.com'"" at ""//tr[td[text()='auto@nesto.com']]/td[2]"""

# TC_CUST_06 - Delete Customer
$ws.Range("A17").Value = "TC_CUST_06"
$ws.Range("B17").Value = "Delete Customer"
$ws.Range("C17").Value = "1.Click ""Delete Button"" at ""//tr[td[text()='auto@nesto.com']]//a[contains(@class, 'btn-outline-danger')]"""
$ws.Range("C18").Value = "2.Verify text ""{DB_QUERY}SELECT COUNT(*) FROM customers WHERE email='auto@nesto.com'"" at ""//table/tbody/tr[td[text()='auto@nesto.com']]"""

# Make Customer_Tests the active sheet/tab, with C19 (just past the new data)
# as the active selection - matches the saved workbook view state.
$ws.Activate()
$ws.Range("C19").Select()
